# Change derivation maximum coefficients and add new method to cost secondary track.
#
# Order of operations matters here: the underlying OOXML shared-string table and
# cellXfs (style) table are append-only logs, so the sequence in which new
# strings/number-formats are introduced determines their final index. We
# replicate the exact order needed to land on the same indices as the target
# workbook. Likewise, the last worksheet that gets Select()/Activate()'d ends
# up as the active tab, so infrastructure is activated last.

$wb = $excel.ActiveWorkbook

$wsMobility       = $wb.Worksheets.Item("mobility")
$wsInfrastructure = $wb.Worksheets.Item("infrastructure")
$wsTime           = $wb.Worksheets.Item("time")
$wsDerivation     = $wb.Worksheets.Item("derivation")
$wsCategories     = $wb.Worksheets.Item("categories")

# --- categories: apply percent (no-decimals) number format to B2:B6 ---------
# This is the first new cellXfs created (numFmtId=9, no alignment) -> style 11.
$wsCategories.Range("B2:B6").NumberFormat = "0%"

# --- infrastructure: new rows + corrected description ----------------------
# String append order: low_quality_track_price, its description, the
# corrected high_quality description (replacing the old typo'd string, which
# gets garbage collected since it becomes unused), gross_main_min_density,
# and its description.
$wsInfrastructure.Range("A16").Value = "low_quality_track_price"
$wsInfrastructure.Range("B16").Value = 200000
$wsInfrastructure.Range("B16").NumberFormat = "#,##0"
$wsInfrastructure.Range("C16").Value = "The price of 1km of low quality track (USD/km)."

$wsInfrastructure.Range("C12").Value = "The price of 1km of high quality track (USD/km)."

$wsInfrastructure.Range("A17").Value = "gross_main_min_density"
$wsInfrastructure.Range("B17").Value = 1333333.3333333333
$wsInfrastructure.Range("B17").NumberFormat = "#,##0"
$wsInfrastructure.Range("C17").Value = "Minimum gross density to consider a link as being a main track (ton-km/ton = ton). Used for secondary track eac cost calculation."

# --- view-state: selections on each sheet -----------------------------------
$wsMobility.Range("C24").Select()
$wsDerivation.Range("A5").Select()
$wsCategories.Range("B3").Select()

# infrastructure becomes the active/selected tab (was "time" before), so
# activate/select it last.
$wsInfrastructure.Range("C5").Select()
